# Split the "Going to try these machine settings: " run into several runs
# that add the machine's name/model "(FLUOstar Omega)" before the colon,
# with a spell-check proofErr bracket around "FLUOstar" (the word Word's
# spell-checker doesn't recognise).
#
#   Going to try these machine settings:
# becomes
#   Going to try these machine settings (FLUOstar Omega):

$d = $word.ActiveDocument

# Locate the paragraph that needs editing via its current (pre-edit) text.
$rng = $d.Content
$found = $rng.Find.Execute("Going to try these machine settings: ")
if (-not $found) {
    throw "Anchor text not found - cannot apply edit."
}

# $rng now spans just the matched text; widen it to the whole paragraph so
# we can rebuild the paragraph's run layout with InsertXML (which replaces
# the full contents of the range it is called on).
$para = $rng.Paragraphs(1).Range

# Capture the part of the paragraph that must stay untouched (everything
# after "...machine settings: ", i.e. the "600nm, settling time..." run).
$tailRange = $d.Range($rng.End, $para.End - 1)
$tailText = $tailRange.Text

function Escape-Xml($s) {
    $s.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

$tailXml = Escape-Xml $tailText

$packageXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="5A4CD9E9" w14:textId="64696970" w:rsidR="003A75CD" w:rsidRDefault="003A75CD" w:rsidP="003A75CD"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>Going to try these machine settings</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>FLUOstar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Omega)</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r w:rsidR="00DF0A1E"><w:t xml:space="preserve">$tailXml</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$para.InsertXML($packageXml)

Write-Output "Applied machine-settings run split + FLUOstar Omega insertion."
